$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.989.08"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.006.66"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.05"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.73"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.695"
$ws.Range("E7").Value = "  +11.44%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.744"
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.171"
$ws.Range("E10").Value = "  -3.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000324"
$ws.Range("E11").Value = "  -6.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.27"
$ws.Range("E12").Value = "  +9.54%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.670.21"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.60"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.013.72"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.95"
$ws.Range("E16").Value = "  -2.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.41"
$ws.Range("E17").Value = "  -3.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.132"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.18"
$ws.Range("E19").Value = "  -3.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.871.12"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "425.67"
$ws.Range("E21").Value = "  -3.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "97.61"
$ws.Range("E22").Value = "  +6.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.48"
$ws.Range("E23").Value = "  -1.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.18"
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.31"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.03"
$ws.Range("E26").Value = "  -11.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.66"
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.87"
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.60"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.59"
$ws.Range("E30").Value = "  +24.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.31"
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.129"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "676.73"
$ws.Range("E33").Value = "  -2.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.95"
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.41"
$ws.Range("E35").Value = "  -4.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "42.61"
$ws.Range("E36").Value = "  +4.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.423"
$ws.Range("E37").Value = "  -4.18%  "
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0825"
$ws.Range("E39").Value = "  -9.43%  "
$ws.Range("E40").Value = "  -4.18%  "
$ws.Range("B41").Value = "Dai"
$ws.Range("C41").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.32"
$ws.Range("E42").Value = "  +7.35%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0485"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.150"
$ws.Range("E45").Value = "  +3.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.56"
$ws.Range("E46").Value = "  +2.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.59"
$ws.Range("E47").Value = "  -10.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.31"
$ws.Range("E48").Value = "  -5.42%  "
$ws.Range("E49").Value = "  -11.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.28"
$ws.Range("E50").Value = "  -3.51%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000267"
$ws.Range("E51").Value = "  -5.71%  "
